# Price/volume refresh for the crypto symbol list (scheduled GitHub Actions update).
# Cells are stored as literal text (e.g. "308.50", "-3.14%"), so force the
# NumberFormat to text ("@") before writing each value -- otherwise Excel would
# coerce "308.50" into the number 308.5 or "-3.14%" into a percentage value.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: D2, E2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '308.50'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '-3.14%'

# Row 3: D3, E3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '37.66'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-5.08%'

# Row 4: D4, E4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.065'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '-1.38%'

# Row 5: D5, E5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.07862'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '-4.15%'

# Row 6: D6, E6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.968'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '-2.33%'

# Row 7: D7, E7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '4.356'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '1.96%'

# Row 8: D8, E8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '8.294'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '0.25%'

# Row 9: D9, E9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.113'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.70%'

# Row 10: D10, E10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9286'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '-0.63%'

# Row 11: D11, E11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.1315'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '-7.14%'

# Row 12: D12, E12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.2064'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '3.50%'

# Row 13: D13, E13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.08816'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-2.79%'

# Row 14: D14, E14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.03434'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '-3.50%'

# Row 15: E15
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '-0.85%'

# Row 16: D16, E16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001390'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.82%'

# Row 17: E17
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-3.43%'

# Row 18: E18
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '1,776.17%'

# Row 19: D19, E19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.587'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '-2.10%'

# Row 20: D20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3473'

# Row 21: E21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '1.61%'

# Row 22: E22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '2.11%'

# Row 23: E23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '1.52%'

# Row 24: D24, E24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04326'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '-0.09%'

# Row 25: D25, E25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001224'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '-0.12%'

# Row 26: D26, E26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004601'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-3.66%'

# Row 27: E27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '3.93%'

# Row 39: D39, E39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02283'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '3.06%'

# Row 40: D40, E40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.05050'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '-3.98%'

# Row 41: D41, E41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.007516'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '0.36%'

# Row 42: D42, E42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.009922'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '0.11%'

# Row 43: D43, E43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1353'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-1.96%'

# Row 44: D44, E44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.001982'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-6.47%'

# Row 45: D45, E45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.008798'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '-10.81%'

# Row 46: D46, E46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00006599'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.58%'

# Row 47: D47, E47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.00000000751'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '0.12%'

# Row 48: D48, E48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.003001'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '8.33%'

# Row 50: D50, E50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.00002103'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.12%'

# Row 51: D51, E51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0002003'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '0.12%'
